# Generate Report for Handoff
# Updates the GUID-named handoff/handback file references (old GUID ->
# new GUID, old hash -> new hash) plus the refreshed handoff timestamps
# on the zh-cn and de-de sheets. Cell values and the matching hyperlink
# display text are kept in sync.

$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink {
    param(
        [object]$Worksheet,
        [string]$CellAddress,
        [string]$NewValue
    )

    $range = $Worksheet.Range($CellAddress)
    $range.Value = $NewValue

    $targetAddr = $range.Address()
    foreach ($hl in $Worksheet.Hyperlinks) {
        if ($hl.Range.Address() -eq $targetAddr) {
            $hl.TextToDisplay = $NewValue
        }
    }
}

$oldGuid = "92d3aa9a-d13a-486a-a85b-3ab9efbbe8d0"
$newGuid = "eb02f3a4-3b8c-4b21-8d75-c1fff3ed7e3b"
$oldHash = "21d1b3ffdc3081d6ef4d67770c77dfa8b1a0e6fd"
$newHash = "d5209ffb99d0038c24af071d96c7e9a21945fa85"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
Set-CellAndHyperlink $wsOverview "A2" ($newGuid + ".md")

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
Set-CellAndHyperlink $wsZh "A2" ($newGuid + ".md")
Set-CellAndHyperlink $wsZh "C2" ($newGuid + "." + $newHash + ".zh-cn.xlf")
Set-CellAndHyperlink $wsZh "D2" "2016-02-29 13:36:30"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
Set-CellAndHyperlink $wsDe "A2" ($newGuid + ".md")
Set-CellAndHyperlink $wsDe "C2" ($newGuid + "." + $newHash + ".de-de.xlf")
Set-CellAndHyperlink $wsDe "D2" "2016-02-29 13:36:39"
